$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.342025122965357
$ws.Range("F2").Value = 4.976228284629165

$ws.Range("C3").Value = 5.091834212275085
$ws.Range("F3").Value = 5.73961504660474

$ws.Range("C4").Value = 5.700275837390827
$ws.Range("F4").Value = 6.492324295674023

$ws.Range("C5").Value = 21
$ws.Range("F5").Value = 34

$ws.Range("C7").Value = 5.141864658058016
$ws.Range("F7").Value = 5.746000290563352

$ws.Range("F8").Value = 4.3

$ws.Range("C9").Value = 1.058498506621991
$ws.Range("F9").Value = 0.8716001340945519
